$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 1, shifting all existing data down by one row.
$ws.Rows.Item(1).Insert()

# Set the new header cell.
$ws.Range("A1").Value = "data"

# Update the view/selection state to match the saved file.
$ws.Range("A2").Select()
